$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Data fix: change the local workstation path used for the
#    "RestoreBy0.30%HF" sheet to the shared server path, matching the
#    other sheets (RestoreBy1.00%HF, RestoreBy0.10%HF, etc.) which
#    already use /scratch/utkur/utkarsh/RestorePillars/exSitu/...
# ------------------------------------------------------------------
$wsTarget = $wb.Worksheets.Item("RestoreBy0.30%HF")
$oldPrefix = "/home/utkarsh/Desktop/exSitu"
$newPrefix = "/scratch/utkur/utkarsh/RestorePillars/exSitu"

for ($row = 2; $row -le 61; $row++) {
    $cell = $wsTarget.Cells.Item($row, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $val.Contains($oldPrefix)) {
        $cell.Value = $val.Replace($oldPrefix, $newPrefix)
    }
}

# ------------------------------------------------------------------
# 2. View-state updates: move the selection/scroll position on each
#    sheet. Each worksheet is activated, and the desired cell is
#    selected so Excel records the new activeCell/sqref (and, for
#    frozen-pane sheets, recalculates the pane's topLeftCell).
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("trainingData")
$ws.Activate()
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("Plasma+RestoreByWater")
$ws.Activate()
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("RestoreByWater")
$ws.Activate()
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("RestoreBy1.00%HF")
$ws.Activate()
$ws.Range("A2").Select()

$ws = $wb.Worksheets.Item("RestoreBy0.10%HF")
$ws.Activate()
$ws.Range("A1").Select()

# RestoreBy0.30%HF keeps its view/selection at A1, and stays the
# active (tab-selected) sheet, matching the workbook's activeTab.
$wsTarget.Activate()
$wsTarget.Range("A1").Select()
